$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 73500
$ws.Range("J87").Value = 73500
$ws.Range("L87").Value = 73500
$ws.Range("N87").Value = -75996

$ws.Range("H90").Value = 73500
$ws.Range("J90").Value = 73500
$ws.Range("L90").Value = 220500
$ws.Range("N90").Value = -232980

$ws.Range("H112").Value = 5620.278
$ws.Range("I112").Value = 2442.5
$ws.Range("J112").Value = 6017.5
$ws.Range("K112").Value = 7327.5
$ws.Range("L112").Value = 18052.5
$ws.Range("M112").Value = -6219.5
$ws.Range("N112").Value = -20268.5

$ws.Range("H123").Value = 60000
$ws.Range("J123").Value = 60000
$ws.Range("L123").Value = 60000
$ws.Range("N123").Value = -69800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H74").Value = 2640.9412
$ws.Range("I74").Value = 2041.3793
$ws.Range("J74").Value = 6118.4
$ws.Range("K74").Value = 2041.3793
$ws.Range("L74").Value = 6118.4
$ws.Range("M74").Value = -1167.3793
$ws.Range("N74").Value = -7866.4

$ws.Range("H77").Value = 2640.9412
$ws.Range("I77").Value = 2041.3793
$ws.Range("J77").Value = 6118.4
$ws.Range("K77").Value = 10206.8965
$ws.Range("L77").Value = 30592
$ws.Range("M77").Value = -5838.896500000001
$ws.Range("N77").Value = -39328

$ws.Range("H92").Value = 36516.668
$ws.Range("J92").Value = 36516.668
$ws.Range("L92").Value = 36516.668
$ws.Range("N92").Value = -41508.668

$ws.Range("H122").Value = 2695.8696
$ws.Range("I122").Value = 2660.8
$ws.Range("J122").Value = 2929.6667
$ws.Range("K122").Value = 7982.400000000001
$ws.Range("L122").Value = 8789.000100000001
$ws.Range("M122").Value = -5532.400000000001
$ws.Range("N122").Value = -13689.0001

$ws.Range("H132").Value = 5781.5293
$ws.Range("I132").Value = 4377.643
$ws.Range("K132").Value = 13132.929
$ws.Range("M132").Value = -10602.929

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 61789.668
$ws.Range("I33").Value = 50000
$ws.Range("J33").Value = 64147.6
$ws.Range("K33").Value = 50000
$ws.Range("L33").Value = 64147.6
$ws.Range("M33").Value = -49664
$ws.Range("N33").Value = -64819.6

$ws.Range("H134").Value = 9427
$ws.Range("I134").Value = 6997.5
$ws.Range("J134").Value = 12666.333
$ws.Range("K134").Value = 20992.5
$ws.Range("L134").Value = 37998.999
$ws.Range("M134").Value = -18457.5
$ws.Range("N134").Value = -43068.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 521.18604
$ws.Range("I7").Value = 499.66666
$ws.Range("J7").Value = 570.8461
$ws.Range("K7").Value = 499.66666
$ws.Range("L7").Value = 570.8461
$ws.Range("M7").Value = -386.66666
$ws.Range("N7").Value = -796.8461

$ws.Range("H16").Value = 1066.6666
$ws.Range("I16").Value = 1066.6666
$ws.Range("K16").Value = 1066.6666
$ws.Range("M16").Value = -779.6666

$ws.Range("H22").Value = 4057.4285
$ws.Range("I22").Value = 1706
$ws.Range("J22").Value = 8290
$ws.Range("K22").Value = 1706
$ws.Range("L22").Value = 8290
$ws.Range("M22").Value = -1356
$ws.Range("N22").Value = -8990

$ws.Range("H31").Value = 10893.667
$ws.Range("I31").Value = 3812.6365
$ws.Range("J31").Value = 16885.309
$ws.Range("K31").Value = 3812.6365
$ws.Range("L31").Value = 16885.309
$ws.Range("M31").Value = -3517.6365
$ws.Range("N31").Value = -17475.309

$ws.Range("H34").Value = 10893.667
$ws.Range("I34").Value = 3812.6365
$ws.Range("J34").Value = 16885.309
$ws.Range("K34").Value = 3812.6365
$ws.Range("L34").Value = 16885.309
$ws.Range("M34").Value = -3610.6365
$ws.Range("N34").Value = -17289.309

$ws.Range("H96").Value = 51625
$ws.Range("J96").Value = 51625
$ws.Range("L96").Value = 51625
$ws.Range("N96").Value = -57117

$ws.Range("H105").Value = 1888.3077
$ws.Range("I105").Value = 2506.125
$ws.Range("K105").Value = 2506.125
$ws.Range("M105").Value = -759.125

$ws.Range("H113").Value = 1066.6666
$ws.Range("I113").Value = 1066.6666
$ws.Range("K113").Value = 1066.6666
$ws.Range("M113").Value = 1103.3334

$ws.Range("H122").Value = 4038.8572
$ws.Range("I122").Value = 3292.2
$ws.Range("J122").Value = 5905.5
$ws.Range("K122").Value = 9876.599999999999
$ws.Range("L122").Value = 17716.5
$ws.Range("M122").Value = -7426.599999999999
$ws.Range("N122").Value = -22616.5

$ws.Range("H135").Value = 93896.5
$ws.Range("J135").Value = 93896.5
$ws.Range("L135").Value = 93896.5
$ws.Range("N135").Value = -104036.5

$ws.Range("H140").Value = 79944.75
$ws.Range("J140").Value = 79944.75
$ws.Range("L140").Value = 79944.75
$ws.Range("N140").Value = -90304.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 534.7143
$ws.Range("J2").Value = 442.42856
$ws.Range("L2").Value = 2654.57136
$ws.Range("N2").Value = -2880.57136

$ws.Range("H5").Value = 868.05884
$ws.Range("I5").Value = 645.9
$ws.Range("J5").Value = 1185.4286
$ws.Range("K5").Value = 1937.7
$ws.Range("L5").Value = 3556.2858
$ws.Range("M5").Value = -1825.7
$ws.Range("N5").Value = -3780.2858

$ws.Range("H129").Value = 11907621
$ws.Range("J129").Value = 15154163
$ws.Range("L129").Value = 45462489
$ws.Range("N129").Value = -45472489

$ws.Range("H135").Value = 868.05884
$ws.Range("I135").Value = 645.9
$ws.Range("J135").Value = 1185.4286
$ws.Range("K135").Value = 5813.099999999999
$ws.Range("L135").Value = 10668.8574
$ws.Range("M135").Value = -3278.099999999999
$ws.Range("N135").Value = -15738.8574

$ws.Range("H136").Value = 3305.6
$ws.Range("I136").Value = 3305.6
$ws.Range("K136").Value = 9916.799999999999
$ws.Range("M136").Value = -4816.799999999999

$ws.Range("H141").Value = 6879.7617
$ws.Range("I141").Value = 5195
$ws.Range("K141").Value = 15585
$ws.Range("M141").Value = -10405

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 28599.8
$ws.Range("J43").Value = 28599.8
$ws.Range("L43").Value = 28599.8
$ws.Range("N43").Value = -28901.8

$ws.Range("H125").Value = 26999.5
$ws.Range("J125").Value = 26999.5
$ws.Range("L125").Value = 26999.5
$ws.Range("N125").Value = -31919.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5039.5
$ws.Range("I22").Value = 5039.5
$ws.Range("K22").Value = 5039.5
$ws.Range("M22").Value = -4744.5

$ws.Range("H27").Value = 5039.5
$ws.Range("I27").Value = 5039.5
$ws.Range("K27").Value = 5039.5
$ws.Range("M27").Value = -4932.5

$ws.Range("H40").Value = 6970.393
$ws.Range("I40").Value = 6398.609
$ws.Range("J40").Value = 9600.6
$ws.Range("K40").Value = 6398.609
$ws.Range("L40").Value = 9600.6
$ws.Range("M40").Value = -6262.609
$ws.Range("N40").Value = -9872.6

$ws.Range("H55").Value = 682.9545000000001
$ws.Range("I55").Value = 212
$ws.Range("J55").Value = 1507.125
$ws.Range("K55").Value = 212
$ws.Range("L55").Value = 1507.125
$ws.Range("M55").Value = -39
$ws.Range("N55").Value = -1853.125

$ws.Range("H61").Value = 4804.5454
$ws.Range("I61").Value = 1835
$ws.Range("K61").Value = 1835
$ws.Range("M61").Value = -1633

$ws.Range("H113").Value = 4804.5454
$ws.Range("I113").Value = 1835
$ws.Range("K113").Value = 1835
$ws.Range("M113").Value = 335

$ws.Range("H122").Value = 4546.6
$ws.Range("I122").Value = 4170.2964
$ws.Range("J122").Value = 7933.3335
$ws.Range("K122").Value = 12510.8892
$ws.Range("L122").Value = 23800.0005
$ws.Range("M122").Value = -10060.8892
$ws.Range("N122").Value = -28700.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 15313.2
$ws.Range("J45").Value = 11333
$ws.Range("L45").Value = 11333
$ws.Range("N45").Value = -12315

$ws.Range("H96").Value = 9297.223
$ws.Range("I96").Value = 4893
$ws.Range("K96").Value = 4893
$ws.Range("M96").Value = -3520
